$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1906354515050167
$ws.Range("C2").Value = 0.5418060200668896
$ws.Range("J2").Value = 0.01672240802675585
$ws.Range("P2").Value = 0.1404682274247492
$ws.Range("S2").Value = 0.1103678929765886
$ws.Range("B3").Value = 0.01204819277108434
$ws.Range("C3").Value = 0.006024096385542169
$ws.Range("J3").Value = 0.06024096385542169
$ws.Range("P3").Value = 0.7289156626506024
$ws.Range("S3").Value = 0.1927710843373494
$ws.Range("J4").Value = 0.09523809523809523
$ws.Range("P4").Value = 0.6904761904761905
$ws.Range("S4").Value = 0.2142857142857143
$ws.Range("B6").Value = 0.04739336492890995
$ws.Range("D6").Value = 0.02369668246445497
$ws.Range("F6").Value = 0.06161137440758294
$ws.Range("J6").Value = 0.2890995260663507
$ws.Range("O6").Value = 0.02369668246445497
$ws.Range("Q6").Value = 0.1516587677725119
$ws.Range("R6").Value = 0.06161137440758294
$ws.Range("S6").Value = 0.3412322274881517
$ws.Range("B7").Value = 0.1111111111111111
$ws.Range("D7").Value = 0.01587301587301587
$ws.Range("F7").Value = 0.07407407407407407
$ws.Range("J7").Value = 0.1164021164021164
$ws.Range("O7").Value = 0.04232804232804233
$ws.Range("Q7").Value = 0.1904761904761905
$ws.Range("R7").Value = 0.0582010582010582
$ws.Range("S7").Value = 0.3915343915343915
$ws.Range("B8").Value = 0.1065573770491803
$ws.Range("D8").Value = 0.01434426229508197
$ws.Range("F8").Value = 0.05532786885245902
$ws.Range("J8").Value = 0.1270491803278689
$ws.Range("O8").Value = 0.0389344262295082
$ws.Range("Q8").Value = 0.2069672131147541
$ws.Range("R8").Value = 0.07377049180327869
$ws.Range("S8").Value = 0.3770491803278688
$ws.Range("B9").Value = 0.08035714285714286
$ws.Range("D9").Value = 0.008928571428571428
$ws.Range("F9").Value = 0.05803571428571429
$ws.Range("J9").Value = 0.06696428571428571
$ws.Range("O9").Value = 0.02232142857142857
$ws.Range("Q9").Value = 0.2008928571428572
$ws.Range("R9").Value = 0.09375
$ws.Range("S9").Value = 0.46875
$ws.Range("B10").Value = 0.1046788263283109
$ws.Range("D10").Value = 0.02061855670103093
$ws.Range("F10").Value = 0.06344171292624901
$ws.Range("J10").Value = 0.1459159397303727
$ws.Range("O10").Value = 0.03092783505154639
$ws.Range("Q10").Value = 0.2117367168913561
$ws.Range("R10").Value = 0.06344171292624901
$ws.Range("S10").Value = 0.359238699444885
$ws.Range("G11").Value = 0.1519434628975265
$ws.Range("J11").Value = 0.0989399293286219
$ws.Range("K11").Value = 0.2084805653710247
$ws.Range("L11").Value = 0.5265017667844523
$ws.Range("S11").Value = 0.01413427561837456
$ws.Range("G12").Value = 0.7483870967741936
$ws.Range("J12").Value = 0.1483870967741935
$ws.Range("K12").Value = 0.006451612903225806
$ws.Range("L12").Value = 0.03870967741935484
$ws.Range("S12").Value = 0.05806451612903226
$ws.Range("G13").Value = 0.7857142857142857
$ws.Range("J13").Value = 0.1904761904761905
$ws.Range("S13").Value = 0.02380952380952381
$ws.Range("G14").Value = 0.5
$ws.Range("S14").Value = 0.5
$ws.Range("F15").Value = 0.01298701298701299
$ws.Range("H15").Value = 0.1385281385281385
$ws.Range("I15").Value = 0.05627705627705628
$ws.Range("J15").Value = 0.2727272727272727
$ws.Range("K15").Value = 0.06060606060606061
$ws.Range("M15").Value = 0.02164502164502164
$ws.Range("N15").Value = 0.004329004329004329
$ws.Range("O15").Value = 0.0303030303030303
$ws.Range("S15").Value = 0.4025974025974026
$ws.Range("F16").Value = 0.01666666666666667
$ws.Range("H16").Value = 0.2333333333333333
$ws.Range("I16").Value = 0.08888888888888889
$ws.Range("J16").Value = 0.4055555555555556
$ws.Range("K16").Value = 0.08888888888888889
$ws.Range("O16").Value = 0.02777777777777778
$ws.Range("S16").Value = 0.1055555555555556
$ws.Range("F17").Value = 0.01467505241090147
$ws.Range("H17").Value = 0.1761006289308176
$ws.Range("I17").Value = 0.0880503144654088
$ws.Range("J17").Value = 0.3983228511530398
$ws.Range("K17").Value = 0.1006289308176101
$ws.Range("M17").Value = 0.01886792452830189
$ws.Range("N17").Value = 0.00419287211740042
$ws.Range("O17").Value = 0.05660377358490566
$ws.Range("S17").Value = 0.1425576519916142
$ws.Range("F18").Value = 0.03726708074534162
$ws.Range("H18").Value = 0.1801242236024845
$ws.Range("I18").Value = 0.1055900621118012
$ws.Range("J18").Value = 0.4161490683229814
$ws.Range("K18").Value = 0.08074534161490683
$ws.Range("M18").Value = 0.02484472049689441
$ws.Range("O18").Value = 0.03726708074534162
$ws.Range("S18").Value = 0.1180124223602484
$ws.Range("F19").Value = 0.01771956856702619
$ws.Range("H19").Value = 0.2326656394453005
$ws.Range("I19").Value = 0.1016949152542373
$ws.Range("J19").Value = 0.3705701078582435
$ws.Range("K19").Value = 0.09553158705701079
$ws.Range("M19").Value = 0.01540832049306626
$ws.Range("O19").Value = 0.0600924499229584
$ws.Range("S19").Value = 0.1063174114021572
